$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$cell = $ws.Range("L2")
$fmt = $cell.NumberFormat
$cell.Value = "58285547"
$cell.NumberFormat = $fmt
